$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header row): add/shift text labels
$ws.Range("B1").Value = "顺序"
$ws.Range("C1").Value = "rrr"
$ws.Range("D1").Value = "职业名字"
$ws.Range("E1").Value = "www"

# Row 2: add/shift text labels
$ws.Range("B2").Value = "type_id"
$ws.Range("C2").Value = "display_name"
$ws.Range("D2").Value = "type"
$ws.Range("E2").Value = "id"

# Row 3: fix existing text labels (D3/E3 unchanged)
$ws.Range("B3").Value = "type_id"
$ws.Range("C3").Value = "display_name"

# Row 5: numeric updates
$ws.Range("B5").Value = 92
$ws.Range("C5").Value = 32

# Row 6: numeric updates
$ws.Range("D6").Value = 424
$ws.Range("E6").Value = 43

# Update the active cell / selection shown when the sheet was last saved
$ws.Range("G10").Select()
